# Update the "Deliveries" worksheet: move the "Fiscalia" district row (previously
# the last row of the table) up to be the first data row, shifting the other
# district rows down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deliveries")

# New data for rows 2-6, columns A-E (Distrito, Precio, Tiempo, PedidoMinimo, EnvioGratis)
$data = @(
    @("Fiscalia",      0,  "1 pm llega el menu", 0,  0),
    @("Trujillo",      8,  "30-45 min",          25, 70),
    @("Victor Larco",  12, "35-45 min",          35, 90),
    @("La Esperanza",  10, "40-50 min",          35, 90),
    @("El Porvenir",   14, "50-60 min",          40, 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

# Update the active selection to match the post-edit state (row 2 selected)
$ws.Activate()
$ws.Range("A2:E2").Select()
